# Natmi following Dr Hou advice
# Re-compute the Sema3a -> Plxna2 LR-pair table for a 3x3 cluster grid
# (ECs / FAPs / sCs on both the sending and target side) instead of the
# previous 3x2-row subset, with refreshed expression/specificity metrics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema3a"
$ws.Range("C2").Value = "Plxna2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.397441
$ws.Range("H2").Value = 4.192323
$ws.Range("I2").Value = 0.6676161521996591
$ws.Range("J2").Value = 0.6676161521996592
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 20.186605
$ws.Range("N2").Value = 60.559815
$ws.Range("O2").Value = 0.6134655823028334
$ws.Range("P2").Value = 0.6134655823028335
$ws.Range("Q2").Value = 28.209589477805
$ws.Range("R2").Value = 253.886305300245
$ws.Range("S2").Value = 0.4095595315639409
$ws.Range("T2").Value = 0.4095595315639411

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema3a"
$ws.Range("C3").Value = "Plxna2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.397441
$ws.Range("H3").Value = 4.192323
$ws.Range("I3").Value = 0.6676161521996591
$ws.Range("J3").Value = 0.6676161521996592
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.772365
$ws.Range("N3").Value = 20.317095
$ws.Range("O3").Value = 0.2058103796201654
$ws.Range("P3").Value = 0.2058103796201654
$ws.Range("Q3").Value = 9.463980517965
$ws.Range("R3").Value = 85.175824661685
$ws.Range("S3").Value = 0.137402333724766
$ws.Range("T3").Value = 0.137402333724766

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema3a"
$ws.Range("C4").Value = "Plxna2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.397441
$ws.Range("H4").Value = 4.192323
$ws.Range("I4").Value = 0.6676161521996591
$ws.Range("J4").Value = 0.6676161521996592
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.946877666666667
$ws.Range("N4").Value = 17.840633
$ws.Range("O4").Value = 0.1807240380770012
$ws.Range("P4").Value = 0.1807240380770012
$ws.Range("Q4").Value = 8.310410673384332
$ws.Range("R4").Value = 74.793696060459
$ws.Range("S4").Value = 0.1206542869109522
$ws.Range("T4").Value = 0.1206542869109522

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Sema3a"
$ws.Range("C5").Value = "Plxna2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.2347746666666667
$ws.Range("H5").Value = 0.7043240000000001
$ws.Range("I5").Value = 0.112161700990566
$ws.Range("J5").Value = 0.112161700990566
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 20.186605
$ws.Range("N5").Value = 60.559815
$ws.Range("O5").Value = 0.6134655823028334
$ws.Range("P5").Value = 0.6134655823028335
$ws.Range("Q5").Value = 4.739303460006667
$ws.Range("R5").Value = 42.65373114006
$ws.Range("S5").Value = 0.06880734321025388
$ws.Range("T5").Value = 0.06880734321025389

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema3a"
$ws.Range("C6").Value = "Plxna2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.2347746666666667
$ws.Range("H6").Value = 0.7043240000000001
$ws.Range("I6").Value = 0.112161700990566
$ws.Range("J6").Value = 0.112161700990566
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.772365
$ws.Range("N6").Value = 20.317095
$ws.Range("O6").Value = 0.2058103796201654
$ws.Range("P6").Value = 0.2058103796201654
$ws.Range("Q6").Value = 1.58997973542
$ws.Range("R6").Value = 14.30981761878
$ws.Range("S6").Value = 0.02308404225971187
$ws.Range("T6").Value = 0.02308404225971188

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema3a"
$ws.Range("C7").Value = "Plxna2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.2347746666666667
$ws.Range("H7").Value = 0.7043240000000001
$ws.Range("I7").Value = 0.112161700990566
$ws.Range("J7").Value = 0.112161700990566
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.946877666666667
$ws.Range("N7").Value = 17.840633
$ws.Range("O7").Value = 0.1807240380770012
$ws.Range("P7").Value = 0.1807240380770012
$ws.Range("Q7").Value = 1.396176221899111
$ws.Range("R7").Value = 12.565585997092
$ws.Range("S7").Value = 0.02027031552060027
$ws.Range("T7").Value = 0.02027031552060028

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Sema3a"
$ws.Range("C8").Value = "Plxna2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.4609646666666667
$ws.Range("H8").Value = 1.382894
$ws.Range("I8").Value = 0.2202221468097748
$ws.Range("J8").Value = 0.2202221468097748
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 20.186605
$ws.Range("N8").Value = 60.559815
$ws.Range("O8").Value = 0.6134655823028334
$ws.Range("P8").Value = 0.6134655823028335
$ws.Range("Q8").Value = 9.305311644956667
$ws.Range("R8").Value = 83.74780480461
$ws.Range("S8").Value = 0.1350987075286386
$ws.Range("T8").Value = 0.1350987075286386

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Sema3a"
$ws.Range("C9").Value = "Plxna2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.4609646666666667
$ws.Range("H9").Value = 1.382894
$ws.Range("I9").Value = 0.2202221468097748
$ws.Range("J9").Value = 0.2202221468097748
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 6.772365
$ws.Range("N9").Value = 20.317095
$ws.Range("O9").Value = 0.2058103796201654
$ws.Range("P9").Value = 0.2058103796201654
$ws.Range("Q9").Value = 3.12182097477
$ws.Range("R9").Value = 28.09638877293
$ws.Range("S9").Value = 0.04532400363568754
$ws.Range("T9").Value = 0.04532400363568755

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Sema3a"
$ws.Range("C10").Value = "Plxna2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.4609646666666667
$ws.Range("H10").Value = 1.382894
$ws.Range("I10").Value = 0.2202221468097748
$ws.Range("J10").Value = 0.2202221468097748
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.946877666666667
$ws.Range("N10").Value = 17.840633
$ws.Range("O10").Value = 0.1807240380770012
$ws.Range("P10").Value = 0.1807240380770012
$ws.Range("Q10").Value = 2.741300481322444
$ws.Range("R10").Value = 24.671704331902
$ws.Range("S10").Value = 0.03979943564544867
$ws.Range("T10").Value = 0.03979943564544869
